# The deck's theme (ppt/theme/theme1.xml) was swapped from the custom
# "Integral" color scheme to the stock PowerPoint "Office Theme" color
# scheme. Reproduce that by rewriting each of the twelve theme colour
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) exposed via the
# classic ColorScheme object, which iron_native maps straight onto
# <a:clrScheme> in ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ColorScheme

function HexRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# 1 dk1      000000
$cs.Colors(1).RGB  = HexRGB 0x00 0x00 0x00
# 2 lt1      FFFFFF
$cs.Colors(2).RGB  = HexRGB 0xFF 0xFF 0xFF
# 3 dk2      44546A
$cs.Colors(3).RGB  = HexRGB 0x44 0x54 0x6A
# 4 lt2      E7E6E6
$cs.Colors(4).RGB  = HexRGB 0xE7 0xE6 0xE6
# 5 accent1  5B9BD5
$cs.Colors(5).RGB  = HexRGB 0x5B 0x9B 0xD5
# 6 accent2  ED7D31
$cs.Colors(6).RGB  = HexRGB 0xED 0x7D 0x31
# 7 accent3  A5A5A5
$cs.Colors(7).RGB  = HexRGB 0xA5 0xA5 0xA5
# 8 accent4  FFC000
$cs.Colors(8).RGB  = HexRGB 0xFF 0xC0 0x00
# 9 accent5  4472C4
$cs.Colors(9).RGB  = HexRGB 0x44 0x72 0xC4
# 10 accent6 70AD47
$cs.Colors(10).RGB = HexRGB 0x70 0xAD 0x47
# 11 hlink   0563C1
$cs.Colors(11).RGB = HexRGB 0x05 0x63 0xC1
# 12 folHlink 954F72
$cs.Colors(12).RGB = HexRGB 0x95 0x4F 0x72
